# Apply cryptos list update (commit: Wed Jun 21 16:27:24 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values are written as literal text (avoids Excel
# auto-converting numeric-looking strings like "1.130" or "0.4160" into numbers)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.114.72"
$ws.Range("E2").Value = "  +10.89%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.871.51"
$ws.Range("E3").Value = "  +7.92%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9949"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.16"
$ws.Range("E5").Value = "  +4.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9945"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4981"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.48"
$ws.Range("E8").Value = "  +7.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2849"
$ws.Range("E9").Value = "  +10.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06549"
$ws.Range("E10").Value = "  +6.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.855.12"
$ws.Range("E11").Value = "  +7.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "17.04"
$ws.Range("E12").Value = "  +6.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07158"
$ws.Range("E13").Value = "  +3.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6694"
$ws.Range("E14").Value = "  +11.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "86.25"
$ws.Range("E15").Value = "  +12.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.772"
$ws.Range("E16").Value = "  +7.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.053.71"
$ws.Range("E17").Value = "  +10.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9953"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007473"
$ws.Range("E19").Value = "  +6.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.55"
$ws.Range("E20").Value = "  +10.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9947"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.089.24"
$ws.Range("E22").Value = "  +7.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.669"
$ws.Range("E23").Value = "  +6.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.501"
$ws.Range("E24").Value = "  +7.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.966"
$ws.Range("E25").Value = "  +7.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "143.84"
$ws.Range("E26").Value = "  +1.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "133.04"
$ws.Range("E27").Value = "  +24.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.79"
$ws.Range("E28").Value = "  +10.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.949"
$ws.Range("E29").Value = "  +7.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.402"
$ws.Range("E30").Value = "  +1.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.237"
$ws.Range("E31").Value = "  +7.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08612"
$ws.Range("E32").Value = "  +8.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.876"
$ws.Range("E33").Value = "  +5.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05026"
$ws.Range("E34").Value = "  +6.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.130"
$ws.Range("E35").Value = "  +11.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6853"
$ws.Range("E36").Value = "  +11.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.685"
$ws.Range("E37").Value = "  +3.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.313"
$ws.Range("E38").Value = "  +15.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.755"
$ws.Range("E39").Value = "  +8.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9504"
$ws.Range("E40").Value = "  +3.07%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01622"
$ws.Range("E41").Value = "  +8.99%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.123"
$ws.Range("E42").Value = "  +7.17%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "103.69"
$ws.Range("E43").Value = "  +4.58%  "
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9954"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4160"
$ws.Range("E45").Value = "  +9.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.414"
$ws.Range("E46").Value = "  +8.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1243"
$ws.Range("E47").Value = "  +8.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05614"
$ws.Range("E48").Value = "  +5.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "32.41"
$ws.Range("E49").Value = "  +8.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.283"
$ws.Range("E50").Value = "  +6.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.335"
$ws.Range("E51").Value = "  +7.81%  "
